# Weekly update: insert two new daily price observations into the
# "Hortaliza, Femacal de La Calera - Arveja Verde" sheet.
#
# The sheet holds one row per market day (rows 2..63, newest edits keep
# getting prepended near the top of the date-sorted block). This commit
# adds two freshly-scraped rows:
#   - a new row at position 7 (pushes the former rows 7..63 down by one)
#   - a new row at position 22 in the *new* numbering (pushes the rest
#     down by one more)
# Excel's Rows.Insert() naturally shifts everything below down and
# copies the row-above formatting (so column D keeps its date style),
# which is exactly what the target file shows: the whole block from the
# old row 7 onward slides down by one or two rows, with dimension
# growing from A1:R63 to A1:R65.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the first new row at row 7 -------------------------------
$ws.Rows.Item(7).Insert()

$ws.Cells.Item(7, 1).Value = 3
$ws.Cells.Item(7, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(7, 3).Value = "Coquimbo"
$ws.Cells.Item(7, 4).Value = 44847
$ws.Cells.Item(7, 5).Value = 5
$ws.Cells.Item(7, 6).Value = 100112022
$ws.Cells.Item(7, 7).Value = "Arveja Verde"
$ws.Cells.Item(7, 8).Value = "Perfection"
$ws.Cells.Item(7, 9).Value = "Primera"
$ws.Cells.Item(7, 10).Value = 71
$ws.Cells.Item(7, 11).Value = 30000
$ws.Cells.Item(7, 12).Value = 31000
$ws.Cells.Item(7, 13).Value = 30493
$ws.Cells.Item(7, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(7, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(7, 16).Value = 1220
$ws.Cells.Item(7, 17).Value = 25
$ws.Cells.Item(7, 18).Value = "Hortaliza"

# --- Insert the second new row at row 22 (post first-insert numbering)
$ws.Rows.Item(22).Insert()

$ws.Cells.Item(22, 1).Value = 3
$ws.Cells.Item(22, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(22, 3).Value = "Coquimbo"
$ws.Cells.Item(22, 4).Value = 44848
$ws.Cells.Item(22, 5).Value = 5
$ws.Cells.Item(22, 6).Value = 100112022
$ws.Cells.Item(22, 7).Value = "Arveja Verde"
$ws.Cells.Item(22, 8).Value = "Perfection"
$ws.Cells.Item(22, 9).Value = "Primera"
$ws.Cells.Item(22, 10).Value = 38
$ws.Cells.Item(22, 11).Value = 30000
$ws.Cells.Item(22, 12).Value = 30000
$ws.Cells.Item(22, 13).Value = 30000
$ws.Cells.Item(22, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(22, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(22, 16).Value = 1200
$ws.Cells.Item(22, 17).Value = 25
$ws.Cells.Item(22, 18).Value = "Hortaliza"

"done"
